$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("alr")

$ws.Range("W2").Value = "R 162,R162,R136,R 147,R147,R 136"
$ws.Range("V3").Value = "IX. PHYSICAL,V. RESIDENT CARE"
$ws.Range("W4").Value = "R145,R128,R136,R179,R200"
$ws.Range("Z4").Value = "F, D"
$ws.Range("W8").Value = "R 266,R100,R266"
$ws.Range("W9").Value = "R213,R190,R126,R179"
$ws.Range("W10").Value = "R999,R 247,R251,R 251,R247,R 145,R145"

$wb.Save()
